# Generate Report for Handback
# Swap the "77232830-..." and "2e932acd-..." file rows on each sheet and
# update their handback status / datetime to reflect the new handback.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "2e932acd-e47f-4f3c-8372-e61745a5bd03.md"
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack

$wsOverview.Range("A3").Value = "77232830-5d71-4781-bf09-c9d381f516af.md"
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "2e932acd-e47f-4f3c-8372-e61745a5bd03.md"
$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Range("D2").Value = "2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.zh-cn.xlf"
$wsZhCn.Range("F2").Value = "2e932acd-e47f-4f3c-8372-e61745a5bd03.md"
$wsZhCn.Range("G2").Value = "2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-03-25 10:55:29"

$wsZhCn.Range("A3").Value = "77232830-5d71-4781-bf09-c9d381f516af.md"
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("D3").Value = "77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.zh-cn.xlf"
$wsZhCn.Range("F3").Value = "77232830-5d71-4781-bf09-c9d381f516af.md"
$wsZhCn.Range("G3").Value = "77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-03-25 10:55:29"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "2e932acd-e47f-4f3c-8372-e61745a5bd03.md"
$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Range("D2").Value = "2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.de-de.xlf"
$wsDeDe.Range("F2").Value = "2e932acd-e47f-4f3c-8372-e61745a5bd03.md"
$wsDeDe.Range("G2").Value = "2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-03-25 10:55:44"

$wsDeDe.Range("A3").Value = "77232830-5d71-4781-bf09-c9d381f516af.md"
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("D3").Value = "77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.de-de.xlf"
$wsDeDe.Range("F3").Value = "77232830-5d71-4781-bf09-c9d381f516af.md"
$wsDeDe.Range("G3").Value = "77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-03-25 10:55:44"
